$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4254.9414
$ws.Range("J43").Value = 4150
$ws.Range("L43").Value = 4150
$ws.Range("N43").Value = -4288

$ws.Range("H92").Value = 805.7143
$ws.Range("J92").Value = 796
$ws.Range("L92").Value = 796
$ws.Range("N92").Value = -3292

$ws.Range("H137").Value = 21919.6
$ws.Range("I137").Value = 9803.25
$ws.Range("J137").Value = 29997.166
$ws.Range("K137").Value = 29409.75
$ws.Range("L137").Value = 89991.49800000001
$ws.Range("M137").Value = -26859.75
$ws.Range("N137").Value = -95091.49800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 10072.111
$ws.Range("I31").Value = 10072.111
$ws.Range("K31").Value = 10072.111
$ws.Range("M31").Value = -9778.111000000001

$ws.Range("H32").Value = 2239.6562
$ws.Range("I32").Value = 1886.4921
$ws.Range("K32").Value = 1886.4921
$ws.Range("M32").Value = -1599.4921

$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26248

$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81240

$ws.Range("H74").Value = 18139.143
$ws.Range("J74").Value = 8002.6
$ws.Range("L74").Value = 8002.6
$ws.Range("N74").Value = -9750.6

$ws.Range("H77").Value = 18139.143
$ws.Range("J77").Value = 8002.6
$ws.Range("L77").Value = 40013
$ws.Range("N77").Value = -48749

$ws.Range("H102").Value = 1510.68
$ws.Range("I102").Value = 1633.3158
$ws.Range("J102").Value = 1122.3334
$ws.Range("K102").Value = 1633.3158
$ws.Range("L102").Value = 1122.3334
$ws.Range("M102").Value = -11.31580000000008
$ws.Range("N102").Value = -4366.3334

$ws.Range("H103").Value = 48766.332
$ws.Range("J103").Value = 48766.332
$ws.Range("L103").Value = 48766.332
$ws.Range("N103").Value = -51110.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3681.2334
$ws.Range("I105").Value = 3293.8462
$ws.Range("K105").Value = 3293.8462
$ws.Range("M105").Value = -1546.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H22").Value = 266.1905
$ws.Range("I22").Value = 270.5263
$ws.Range("K22").Value = 270.5263
$ws.Range("M22").Value = 79.47370000000001

$ws.Range("H31").Value = 3364.389
$ws.Range("I31").Value = 1691.75
$ws.Range("J31").Value = 3842.2856
$ws.Range("K31").Value = 1691.75
$ws.Range("L31").Value = 3842.2856
$ws.Range("M31").Value = -1396.75
$ws.Range("N31").Value = -4432.2856

$ws.Range("H34").Value = 3364.389
$ws.Range("I34").Value = 1691.75
$ws.Range("J34").Value = 3842.2856
$ws.Range("K34").Value = 1691.75
$ws.Range("L34").Value = 3842.2856
$ws.Range("M34").Value = -1489.75
$ws.Range("N34").Value = -4246.2856

$ws.Range("H58").Value = 2986.5088
$ws.Range("J58").Value = 3742.3667
$ws.Range("L58").Value = 3742.3667
$ws.Range("N58").Value = -4148.3667

$ws.Range("H94").Value = 7957.25
$ws.Range("I94").Value = 7214.8335
$ws.Range("J94").Value = 8699.666999999999
$ws.Range("K94").Value = 7214.8335
$ws.Range("L94").Value = 8699.666999999999
$ws.Range("M94").Value = -6763.8335
$ws.Range("N94").Value = -9601.666999999999

$ws.Range("H136").Value = 2986.5088
$ws.Range("J136").Value = 3742.3667
$ws.Range("L136").Value = 11227.1001
$ws.Range("N136").Value = -16327.1001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29917606
$ws.Range("J4").Value = 200046.4
$ws.Range("L4").Value = 600139.2
$ws.Range("N4").Value = -600363.2

$ws.Range("H5").Value = 1876.6471
$ws.Range("J5").Value = 2038.4667
$ws.Range("L5").Value = 6115.4001
$ws.Range("N5").Value = -6339.4001

$ws.Range("H37").Value = 113329.664
$ws.Range("J37").Value = 113329.664
$ws.Range("L37").Value = 339988.992
$ws.Range("N37").Value = -340212.992

$ws.Range("H39").Value = 7495.8
$ws.Range("J39").Value = 8992
$ws.Range("L39").Value = 26976
$ws.Range("N39").Value = -27564

$ws.Range("H40").Value = 48.466667
$ws.Range("I40").Value = 45.272728
$ws.Range("J40").Value = 57.25
$ws.Range("K40").Value = 181.090912
$ws.Range("L40").Value = 229
$ws.Range("M40").Value = -112.090912
$ws.Range("N40").Value = -367

$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H42").Value = 300
$ws.Range("I42").Value = 300
$ws.Range("K42").Value = 900
$ws.Range("M42").Value = -366

$ws.Range("H44").Value = 976
$ws.Range("J44").Value = 976
$ws.Range("L44").Value = 2928
$ws.Range("N44").Value = -3724

$ws.Range("H46").Value = 544.1667
$ws.Range("I46").Value = 466.25
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 1398.75
$ws.Range("L46").Value = 2100
$ws.Range("M46").Value = -1307.75
$ws.Range("N46").Value = -2282

$ws.Range("H86").Value = 490.42856
$ws.Range("J86").Value = 490.42856
$ws.Range("L86").Value = 1471.28568
$ws.Range("N86").Value = -3843.28568

$ws.Range("H89").Value = 490.42856
$ws.Range("J89").Value = 490.42856
$ws.Range("L89").Value = 4413.85704
$ws.Range("N89").Value = -16269.85704

$ws.Range("H135").Value = 1876.6471
$ws.Range("J135").Value = 2038.4667
$ws.Range("L135").Value = 18346.2003
$ws.Range("N135").Value = -23416.2003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2129.3157
$ws.Range("I82").Value = 1741.2858
$ws.Range("J82").Value = 2355.6667
$ws.Range("K82").Value = 1741.2858
$ws.Range("L82").Value = 2355.6667
$ws.Range("M82").Value = -1380.2858
$ws.Range("N82").Value = -3077.6667

$ws.Range("H85").Value = 2129.3157
$ws.Range("I85").Value = 1741.2858
$ws.Range("J85").Value = 2355.6667
$ws.Range("K85").Value = 1741.2858
$ws.Range("L85").Value = 2355.6667
$ws.Range("M85").Value = -493.2858000000001
$ws.Range("N85").Value = -4851.6667

$ws.Range("H101").Value = 21299.8
$ws.Range("J101").Value = 21299.8
$ws.Range("L101").Value = 21299.8
$ws.Range("N101").Value = -27789.8

$ws.Range("H132").Value = 5440.028
$ws.Range("I132").Value = 5373.3213
$ws.Range("K132").Value = 16119.9639
$ws.Range("M132").Value = -13589.9639

$ws.Range("H136").Value = 4961.5557
$ws.Range("I136").Value = 4237.12
$ws.Range("K136").Value = 12711.36
$ws.Range("M136").Value = -10161.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 68818.336
$ws.Range("I64").Value = 66457
$ws.Range("K64").Value = 66457
$ws.Range("M64").Value = -66209

$ws.Range("H67").Value = 68818.336
$ws.Range("I67").Value = 66457
$ws.Range("K67").Value = 66457
$ws.Range("M67").Value = -65599

$ws.Range("H103").Value = 28266.666
$ws.Range("J103").Value = 28266.666
$ws.Range("L103").Value = 28266.666
$ws.Range("N103").Value = -30610.666
